$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.285.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.756.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.83%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.58"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.750.29"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.87%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.42%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.39"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.483"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.19"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000257"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.376.53"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.754.84"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.291.47"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.56"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.57"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "505.35"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.23"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.724"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.62"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.60"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.66%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.09"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.34"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000134"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.45%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.50"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.56%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.94"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.66"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.67%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.16"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.19%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.11%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.10"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +14.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.08"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "50.01"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "45.53"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "433.62"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.64"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.962.33"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0364"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.52"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.04%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.42"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.49"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.66%  "
